# fix: fix allure bug
$wb = $excel.ActiveWorkbook
$usersSheet = $wb.Worksheets.Item("users")

# --- Workbook window view size change (bookViews/workbookView) ---
# windowWidth 28800 -> 25600, windowHeight 12280 -> 10480
$win = $wb.Windows.Item(1)
$win.Width = 25600
$win.Height = 10480
$excel.ActiveWindow.Width = 25600
$excel.ActiveWindow.Height = 10480

# --- Text content fixes (shared strings) ---
# "block" -> "blocker" (column E / severity) for every data row
$usedRange = $usersSheet.UsedRange
$rowCount = $usedRange.Rows.Count
for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $usersSheet.Cells.Item($r, 5) # column E
    if ($cell.Value2 -eq "block") {
        $cell.Value2 = "blocker"
    }
}

# D9 "修改用户" -> "查找用户"
$usersSheet.Cells.Item(9, 4).Value2 = "查找用户"

# --- Cell style fix: E4 and E6 switch from style index 5 (blue fill) to style index 4 (yellow fill) ---
$usersSheet.Range("E4").Interior.Color = 65535
$usersSheet.Range("E6").Interior.Color = 65535

# --- Selection change: E13 -> F4 ---
$usersSheet.Activate()
$usersSheet.Range("F4").Select()
